# GDD-New.docx -- "Updated the Greenfoot world"
#
# Renames the enemy from "ghost(s)" to "snake(s)" throughout the document,
# drops the "(Similar to Pac-Man)" aside, swaps the first/third level
# themes (Cyber <-> Desert), and tidies a couple of other small wording
# bits, all per the authoritative diff.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    # MatchCase on, whole string (not whole-word) match, replace all
    # occurrences found in the main content story.
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1) Bold "ghosts" -> bold "snakes " (note trailing space; the following
#    "(Similar to Pac-Man) " aside is dropped right after, see step 2).
#    Isolate this as its own whole-word match so the bold run keeps its
#    formatting and doesn't bleed onto the plain-text "which are trying…"
#    that follows it.
$d.Content.Find.Execute("ghosts", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "snakes ", 1) | Out-Null

# 2) Drop the "(Similar to Pac-Man) " aside that used to follow "snakes ".
#    (There was a lone space-only run between "ghosts"/"snakes " and the
#    aside, so the find text below needs the leading space to consume it.)
Replace-Text " (Similar to Pac-Man) which" "which"

# 3) "without getting caught by the ghosts or without running out of time."
Replace-Text "without getting caught by the ghosts or without running out of time." "without getting caught by the snakes or without running out of time."

# 4) First level: ghosts -> snakes, world theme cyber -> desert.
Replace-Text ", a small number of ghosts and matches in a cyber world. For the " ", a small number of snakes and matches in a desert world. For the "

# 5) Second level: ghosts -> snakes (both mentions), theme unchanged.
Replace-Text ", An increased number of ghosts and matches, the speed of the ghosts increases, in a grassy location. For the " ", An increased number of snakes and matches, the speed of the snakes increases, in a grassy location. For the "

# 6) Third level: ghosts/ghost -> snakes, theme desert -> cyber.
Replace-Text ", A much more increased number of ghosts and matches, the ghost speed has been increased even more, in a desert area, " ", A much more increased number of snakes and matches, the snakes speed has been increased even more, in a cyber area, "

# 7) Win/Lose condition paragraph.
Replace-Text "gets caught by one of the ghosts. " "gets caught by one of the snakes."

# 8) Enemies bullet.
Replace-Text "The ghosts which are roaming around the map" "The snakes which are roaming around the map"

# 9) Goals of the character bullet.
Replace-Text "Do not get caught by any ghosts roaming around the map" "Do not get caught by any snakes roaming around the map"

# 10) Losing lives bullet.
Replace-Text "they touch a ghost." "they touch a snake."

# 11) Developer roles bullet.
Replace-Text "Character and Ghosts)" "Character and snakes)"

# 12) Levels: swap the first and third level themes.
$d.Content.Find.Execute("Cyber", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "~~TmpCyber~~", 1) | Out-Null
$d.Content.Find.Execute("Desert", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "Cyber", 1) | Out-Null
$d.Content.Find.Execute("~~TmpCyber~~", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "Desert", 1) | Out-Null

# 13) Random Motion bullet.
Replace-Text "Ghosts are moving in random motion" "Snakes are moving in random motion"

# 14) Animation bullet.
Replace-Text "or the ghosts (could be both)" "or the snakes (could be both)"

# 15) Appearing bullet.
Replace-Text "an additional ghost will randomly appear on the map" "an additional snake will randomly appear on the map"
